$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 673.5172858837635

$ws.Range("B3").Value = 0.127881588408715
$ws.Range("C3").Value = 0.04240448674262143
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 9.63600819768431
